$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name (B) and Link (C) text cells for reordered rows
$cell = $ws.Range('D2')
$cell.NumberFormat = "@"
$cell.Value = '25.572.45'
$cell.Style = "Normal"
$ws.Range('E2').Value = '  -4.31%  '

$cell = $ws.Range('D3')
$cell.NumberFormat = "@"
$cell.Value = '1.808.65'
$cell.Style = "Normal"
$ws.Range('E3').Value = '  -3.14%  '

$cell = $ws.Range('D4')
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range('E4').Value = '  +0.08%  '

$cell = $ws.Range('D5')
$cell.NumberFormat = "@"
$cell.Value = '275.77'
$cell.Style = "Normal"
$ws.Range('E5').Value = '  -8.28%  '

$cell = $ws.Range('D6')
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range('E6').Value = '  +0.04%  '

$cell = $ws.Range('D7')
$cell.NumberFormat = "@"
$cell.Value = '0.5000'
$cell.Style = "Normal"
$ws.Range('E7').Value = '  -5.89%  '

$cell = $ws.Range('D8')
$cell.NumberFormat = "@"
$cell.Value = '0.3485'
$cell.Style = "Normal"
$ws.Range('E8').Value = '  -6.55%  '

$cell = $ws.Range('D9')
$cell.NumberFormat = "@"
$cell.Value = '44.24'
$cell.Style = "Normal"
$ws.Range('E9').Value = '  -2.58%  '

$cell = $ws.Range('D10')
$cell.NumberFormat = "@"
$cell.Value = '0.06636'
$cell.Style = "Normal"
$ws.Range('E10').Value = '  -7.46%  '

$cell = $ws.Range('D11')
$cell.NumberFormat = "@"
$cell.Value = '19.48'
$cell.Style = "Normal"
$ws.Range('E11').Value = '  -9.25%  '

$cell = $ws.Range('D12')
$cell.NumberFormat = "@"
$cell.Value = '0.8031'
$cell.Style = "Normal"
$ws.Range('E12').Value = '  -9.61%  '

$cell = $ws.Range('D13')
$cell.NumberFormat = "@"
$cell.Value = '0.07881'
$cell.Style = "Normal"
$ws.Range('E13').Value = '  -3.67%  '

$cell = $ws.Range('D14')
$cell.NumberFormat = "@"
$cell.Value = '1.803.05'
$cell.Style = "Normal"
$ws.Range('E14').Value = '  -3.80%  '

$cell = $ws.Range('D15')
$cell.NumberFormat = "@"
$cell.Value = '5.029'
$cell.Style = "Normal"
$ws.Range('E15').Value = '  -5.04%  '

$cell = $ws.Range('D16')
$cell.NumberFormat = "@"
$cell.Value = '86.54'
$cell.Style = "Normal"
$ws.Range('E16').Value = '  -6.55%  '

$cell = $ws.Range('D17')
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Range('E17').Value = '  +0.10%  '

$cell = $ws.Range('D18')
$cell.NumberFormat = "@"
$cell.Value = '13.93'
$cell.Style = "Normal"
$ws.Range('E18').Value = '  -5.75%  '

$cell = $ws.Range('D19')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range('E19').Value = '  +0.00%  '

$cell = $ws.Range('D20')
$cell.NumberFormat = "@"
$cell.Value = '0.000007925'
$cell.Style = "Normal"
$ws.Range('E20').Value = '  -6.64%  '

$cell = $ws.Range('D21')
$cell.NumberFormat = "@"
$cell.Value = '25.637.80'
$cell.Style = "Normal"
$ws.Range('E21').Value = '  -4.24%  '

$cell = $ws.Range('D22')
$cell.NumberFormat = "@"
$cell.Value = '4.710'
$cell.Style = "Normal"
$ws.Range('E22').Value = '  -5.34%  '

$cell = $ws.Range('D23')
$cell.NumberFormat = "@"
$cell.Value = '9.888'
$cell.Style = "Normal"
$ws.Range('E23').Value = '  -6.82%  '

$cell = $ws.Range('D24')
$cell.NumberFormat = "@"
$cell.Value = '6.085'
$cell.Style = "Normal"
$ws.Range('E24').Value = '  -4.28%  '

$cell = $ws.Range('D25')
$cell.NumberFormat = "@"
$cell.Value = '2.240'
$cell.Style = "Normal"
$ws.Range('E25').Value = '  -2.50%  '

$cell = $ws.Range('D26')
$cell.NumberFormat = "@"
$cell.Value = '142.62'
$cell.Style = "Normal"
$ws.Range('E26').Value = '  -2.16%  '

$cell = $ws.Range('D27')
$cell.NumberFormat = "@"
$cell.Value = '1.658'
$cell.Style = "Normal"
$ws.Range('E27').Value = '  -4.23%  '

$cell = $ws.Range('D28')
$cell.NumberFormat = "@"
$cell.Value = '17.05'
$cell.Style = "Normal"
$ws.Range('E28').Value = '  -5.41%  '

$cell = $ws.Range('D29')
$cell.NumberFormat = "@"
$cell.Value = '108.18'
$cell.Style = "Normal"
$ws.Range('E29').Value = '  -4.73%  '

$cell = $ws.Range('D30')
$cell.NumberFormat = "@"
$cell.Value = '4.253'
$cell.Style = "Normal"
$ws.Range('E30').Value = '  -9.18%  '

$cell = $ws.Range('D31')
$cell.NumberFormat = "@"
$cell.Value = '4.194'
$cell.Style = "Normal"
$ws.Range('E31').Value = '  -9.55%  '

$cell = $ws.Range('D32')
$cell.NumberFormat = "@"
$cell.Value = '0.08728'
$cell.Style = "Normal"
$ws.Range('E32').Value = '  -4.31%  '

$cell = $ws.Range('D33')
$cell.NumberFormat = "@"
$cell.Value = '0.04798'
$cell.Style = "Normal"
$ws.Range('E33').Value = '  -4.34%  '

$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$cell = $ws.Range('D34')
$cell.NumberFormat = "@"
$cell.Value = '2.875'
$cell.Style = "Normal"
$ws.Range('E34').Value = '  -2.49%  '

$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D35')
$cell.NumberFormat = "@"
$cell.Value = '1.122'
$cell.Style = "Normal"
$ws.Range('E35').Value = '  -4.39%  '

$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D36')
$cell.NumberFormat = "@"
$cell.Value = '0.7136'
$cell.Style = "Normal"
$ws.Range('E36').Value = '  -10.75%  '

$cell = $ws.Range('D37')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range('E37').Value = '  -0.09%  '

$cell = $ws.Range('D38')
$cell.NumberFormat = "@"
$cell.Value = '3.124'
$cell.Style = "Normal"
$ws.Range('E38').Value = '  -1.45%  '

$cell = $ws.Range('D39')
$cell.NumberFormat = "@"
$cell.Value = '2.297'
$cell.Style = "Normal"
$ws.Range('E39').Value = '  -14.11%  '

$cell = $ws.Range('D40')
$cell.NumberFormat = "@"
$cell.Value = '0.01829'
$cell.Style = "Normal"
$ws.Range('E40').Value = '  -5.78%  '

$cell = $ws.Range('D41')
$cell.NumberFormat = "@"
$cell.Value = '0.5026'
$cell.Style = "Normal"
$ws.Range('E41').Value = '  -17.02%  '

$cell = $ws.Range('D42')
$cell.NumberFormat = "@"
$cell.Value = '0.9456'
$cell.Style = "Normal"
$ws.Range('E42').Value = '  -11.14%  '

$cell = $ws.Range('D43')
$cell.NumberFormat = "@"
$cell.Value = '114.55'
$cell.Style = "Normal"
$ws.Range('E43').Value = '  +0.10%  '

$cell = $ws.Range('D44')
$cell.NumberFormat = "@"
$cell.Value = '6.155'
$cell.Style = "Normal"
$ws.Range('E44').Value = '  -5.26%  '

$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Range('D45')
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Range('E45').Value = '  +0.01%  '

$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range('D46')
$cell.NumberFormat = "@"
$cell.Value = '7.813'
$cell.Style = "Normal"
$ws.Range('E46').Value = '  -10.22%  '

$cell = $ws.Range('D47')
$cell.NumberFormat = "@"
$cell.Value = '0.1347'
$cell.Style = "Normal"
$ws.Range('E47').Value = '  -9.55%  '

$cell = $ws.Range('D48')
$cell.NumberFormat = "@"
$cell.Value = '0.4390'
$cell.Style = "Normal"
$ws.Range('E48').Value = '  -15.73%  '

$cell = $ws.Range('D49')
$cell.NumberFormat = "@"
$cell.Value = '36.16'
$cell.Style = "Normal"
$ws.Range('E49').Value = '  -3.52%  '

$cell = $ws.Range('D50')
$cell.NumberFormat = "@"
$cell.Value = '9.171'
$cell.Style = "Normal"
$ws.Range('E50').Value = '  -7.53%  '

$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$cell = $ws.Range('D51')
$cell.NumberFormat = "@"
$cell.Value = '0.05823'
$cell.Style = "Normal"
$ws.Range('E51').Value = '  -3.75%  '
